$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.934.31'
$ws.Range('E2').Value = '  +1.84%  '
$ws.Range('D3').Value = '2.523.75'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '594.63'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '175.25'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('D9').Value = '2.520.39'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.150'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +7.59%  '
$ws.Range('E11').Value = '  -0.92%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.98'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.339'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').Value = '2.990.49'
$ws.Range('E14').Value = '  +2.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '25.95'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '68.737.75'
$ws.Range('E16').Value = '  +1.68%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000173'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.11%  '
$ws.Range('D18').Value = '2.518.69'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '364.08'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.47%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.54'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.89%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.98'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('E22').Value = '  +1.74%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '70.55'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.19'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('E26').Value = '  -2.18%  '
$ws.Range('E27').Value = '  -5.80%  '
$ws.Range('D28').Value = '2.646.07'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '513.52'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('D31').Value = '0.0₃0884'
$ws.Range('E31').Value = '  -2.72%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.77'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.22%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.24'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.78'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '162.35'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('E37').Value = '  -3.19%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.68'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.32'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.325'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.53%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.81'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.85%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.37'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.61%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '151.29'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.09%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.58'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.00%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.516'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0741'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.58'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.578'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.19%  '
